$d = $word.ActiveDocument

# Helper: replace the text of a paragraph's trailing run while preserving any
# preceding empty run(s) and the paragraph's own pPr. Find.Execute / Range.Text
# in this runtime tends to merge/drop an immediately-preceding empty <w:r/>,
# so we build a minimal OOXML package fragment and use Range.InsertXML, which
# replaces just the targeted text range's contents (run-for-run) instead of
# the whole paragraph.
function Replace-RunText($oldText, $newText) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        return $false
    }
    $rng = $d.Range($idx, $idx + $oldText.Length)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
    return $true
}

# 1. Title heading ("Play Firekick! ... Slot Game") -- appears twice
#    (Heading1 title near the top, and a bold run near the bottom). Neither
#    occurrence is preceded by an empty run, so a plain Find/Replace across
#    the whole document is safe and replaces both occurrences in one call.
$d.Content.Find.Execute(
    "Play Firekick! Multimax Free | Review of Yggdrasil Gaming's Slot Game",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Firekick! Multimax Slot for Free", 2)

# 2. "What we like" bullet list - reordered and reworded
Replace-RunText "Exciting special features" "Dynamic atmosphere of a stadium on the reels"
Replace-RunText "High maximum possible win" "Exciting special features to enhance the gaming experience"
Replace-RunText "Available on multiple gaming sites" "Possibility of winning up to 10,000 times the bet"
Replace-RunText "Demo mode without registration or download" "Available to play on various online gaming sites"

# 3. "What we don't like" bullet list - reworded
Replace-RunText "Betting limits may not appeal to all players" "Betting limits may not suit all players"
Replace-RunText "Theme may not be of interest to all players" "Limited availability on some gaming sites"

# 4. Italic summary text near the end
$d.Content.Find.Execute(
    "Try Firekick! Multimax for free or real money. Review of Yggdrasil Gaming's online slot game with special features and high maximum possible win.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Firekick! Multimax and play this exciting slot game for free.", 2)
